$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The records in data rows 2-6 get cyclically reshuffled across rows:
#   new row2 = old row5, new row3 = old row4, new row4 = old row2,
#   new row5 = old row6, new row6 = old row3
# Only the columns that actually differ between the shuffled records are
# written here (A, B, E, F, G, H, Q, R, and for rows 5/6 also Y, AA which
# hold plain text dates) -- every other cell in the sheet is left exactly
# as it was.

$ws.Range("A2").Value2 = 111363017
$ws.Range("B2").Value2 = 78578
$ws.Range("E2").Value2 = 6458
$ws.Range("F2").Value2 = "Lunglav"
$ws.Range("G2").Value2 = "Lobaria pulmonaria"
$ws.Range("H2").Value2 = "(L.) Hoffm."
$ws.Range("Q2").Value2 = 593472.3298762256
$ws.Range("R2").Value2 = 6986898.025413335

$ws.Range("A3").Value2 = 111363019
$ws.Range("B3").Value2 = 77268
$ws.Range("E3").Value2 = 228912
$ws.Range("F3").Value2 = "Mörk kolflarnlav"
$ws.Range("G3").Value2 = "Carbonicola myrmecina"
$ws.Range("H3").Value2 = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q3").Value2 = 593439.5664869671
$ws.Range("R3").Value2 = 6986881.627536911

$ws.Range("A4").Value2 = 111363016
$ws.Range("B4").Value2 = 77515
$ws.Range("E4").Value2 = 6425
$ws.Range("F4").Value2 = "Garnlav"
$ws.Range("G4").Value2 = "Alectoria sarmentosa"
$ws.Range("H4").Value2 = "(Ach.) Ach."
$ws.Range("Q4").Value2 = 593439.5664869671
$ws.Range("R4").Value2 = 6986881.627536911

$ws.Range("A5").Value2 = 111363032
$ws.Range("B5").Value2 = 89405
$ws.Range("E5").Value2 = 1202
$ws.Range("F5").Value2 = "Ullticka"
$ws.Range("G5").Value2 = "Phellinidium ferrugineofuscum"
$ws.Range("H5").Value2 = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q5").Value2 = 593471.9433232083
$ws.Range("R5").Value2 = 6986863.916970093
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value2 = "2023-08-05"
$ws.Range("Y5").NumberFormat = "General"
$ws.Range("Y5").Style = "Normal"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value2 = "2023-08-05"
$ws.Range("AA5").NumberFormat = "General"
$ws.Range("AA5").Style = "Normal"

$ws.Range("A6").Value2 = 111363018
$ws.Range("B6").Value2 = 89405
$ws.Range("E6").Value2 = 1202
$ws.Range("F6").Value2 = "Ullticka"
$ws.Range("G6").Value2 = "Phellinidium ferrugineofuscum"
$ws.Range("H6").Value2 = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q6").Value2 = 593479.5069047968
$ws.Range("R6").Value2 = 6986870.044355935
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value2 = "2023-08-06"
$ws.Range("Y6").NumberFormat = "General"
$ws.Range("Y6").Style = "Normal"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value2 = "2023-08-06"
$ws.Range("AA6").NumberFormat = "General"
$ws.Range("AA6").Style = "Normal"
